$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.527.93"
$ws.Range("E2").Value = "  -1.02%  "

$ws.Range("D3").Value = "2.040.72"
$ws.Range("E3").Value = "  -0.57%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'244.48"
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("D6").Value = "'0.661"
$ws.Range("E6").Value = "  +1.24%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").Value = "'53.82"
$ws.Range("E8").Value = "  -6.07%  "

$ws.Range("D9").Value = "'62.86"
$ws.Range("E9").Value = "  +6.36%  "

$ws.Range("D10").Value = "'0.363"
$ws.Range("E10").Value = "  -1.59%  "

$ws.Range("D11").Value = "'0.0743"
$ws.Range("E11").Value = "  -4.24%  "

$ws.Range("E12").Value = "  -3.85%  "

$ws.Range("D13").Value = "'0.938"
$ws.Range("E13").Value = "  +7.61%  "

$ws.Range("D14").Value = "'14.40"
$ws.Range("E14").Value = "  -4.55%  "

$ws.Range("D15").Value = "2.332.97"
$ws.Range("E15").Value = "  -0.76%  "

$ws.Range("D16").Value = "'5.39"
$ws.Range("E16").Value = "  -3.11%  "

$ws.Range("D17").Value = "2.027.57"
$ws.Range("E17").Value = "  -2.67%  "

$ws.Range("D18").Value = "36.381.90"
$ws.Range("E18").Value = "  -1.26%  "

$ws.Range("D19").Value = "'17.02"
$ws.Range("E19").Value = "  -1.77%  "

$ws.Range("D20").Value = "'71.11"
$ws.Range("E20").Value = "  -2.66%  "

$ws.Range("D21").Value = "0.0₃0851"
$ws.Range("E21").Value = "  -4.04%  "

$ws.Range("D22").Value = "'237.62"
$ws.Range("E22").Value = "  +0.71%  "

$ws.Range("D23").Value = "'5.17"
$ws.Range("E23").Value = "  -4.21%  "

$ws.Range("E24").Value = "  +0.14%  "

$ws.Range("D25").Value = "'2.37"
$ws.Range("E25").Value = "  -3.50%  "

$ws.Range("D26").Value = "'2.22"
$ws.Range("E26").Value = "  +1.44%  "

$ws.Range("D27").Value = "'9.17"
$ws.Range("E27").Value = "  -8.96%  "

$ws.Range("D28").Value = "'163.77"
$ws.Range("E28").Value = "  -2.90%  "

$ws.Range("D29").Value = "'19.94"
$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("D30").Value = "'0.121"
$ws.Range("E30").Value = "  -2.56%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.97"
$ws.Range("E31").Value = "  -9.09%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'1.16"
$ws.Range("E32").Value = "  +4.13%  "

$ws.Range("D33").Value = "'0.0596"
$ws.Range("E33").Value = "  -2.55%  "

$ws.Range("D34").Value = "'4.41"
$ws.Range("E34").Value = "  -8.32%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").Value = "'0.0857"
$ws.Range("E36").Value = "  +0.97%  "

$ws.Range("E37").Value = "  -1.55%  "

$ws.Range("D38").Value = "'2.21"
$ws.Range("E38").Value = "  -5.18%  "

$ws.Range("D39").Value = "'1.23"
$ws.Range("E39").Value = "  -6.07%  "

$ws.Range("D40").Value = "'4.89"
$ws.Range("E40").Value = "  +0.71%  "

$ws.Range("D41").Value = "'2.86"
$ws.Range("E41").Value = "  -3.29%  "

$ws.Range("D42").Value = "'0.0213"
$ws.Range("E42").Value = "  -4.19%  "

$ws.Range("D43").Value = "'1.10"
$ws.Range("E43").Value = "  -3.98%  "

$ws.Range("D44").Value = "'93.48"
$ws.Range("E44").Value = "  -3.44%  "

$ws.Range("D45").Value = "'0.0899"
$ws.Range("E45").Value = "  -5.85%  "

$ws.Range("D46").Value = "1.374.15"
$ws.Range("E46").Value = "  +5.40%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'15.67"
$ws.Range("E47").Value = "  -5.69%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'7.43"
$ws.Range("E48").Value = "  +10.47%  "

$ws.Range("D50").Value = "'2.25"
$ws.Range("E50").Value = "  -4.18%  "

$ws.Range("D51").Value = "2.222.20"
$ws.Range("E51").Value = "  -0.66%  "
